# Correct status name:
#  - statut_label "bleu" -> "noir"
#  - statut_name texts reworded from "... et / ou publication posté ..." to "... postés ou publiés ..."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange

# statut_label column: bleu -> noir
$used.Replace("bleu", "noir")

# statut_name column: reworded status descriptions
# (ordered so whole-cell replace matches are unambiguous; Replace matches full cell text)
$used.Replace("pas de résultat ni de publication", "pas de résultat postés ni publiés")
$used.Replace("résultat et / ou publication posté dans les 12 mois", "résultat postés ou publiés dans les 12 mois")
$used.Replace("résultat et / ou publication posté dans les 36 mois", "résultat postés ou publiés dans les 36 mois")
$used.Replace("résultat et / ou publication posté", "résultat postés ou publiés")
